# Update profit.py after running on 2025-08-28
#
# profit.py was re-run for 2025-08-28: it appends the new day's profit row
# to the Sheet1 log and refreshes the Sheet2 "latest day" stats row (date +
# the two computed ratios) to point at that same day.

$wb = $excel.ActiveWorkbook

# --- Sheet1: append the new day's profit row -----------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# The date column stores plain text like "08/27/2025" (see the other rows),
# not a real date value, so force Text formatting before writing the value -
# otherwise Excel auto-converts the "MM/DD/YYYY"-looking string into a date
# serial number. Resetting the style back to "Normal" afterwards keeps the
# cell's value as text without leaving a stray number-format behind.
$ws1.Range("A11").NumberFormat = "@"
$ws1.Cells.Item(11, 1).Value = "08/28/2025"
$ws1.Cells.Item(11, 1).Style = "Normal"
$ws1.Cells.Item(11, 2).Value = 15832.36

# --- Sheet2: refresh the summary row with the new date/stats -------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").NumberFormat = "@"
$ws2.Cells.Item(1, 1).Value = "08/28/2025"
$ws2.Cells.Item(1, 1).Style = "Normal"
$ws2.Cells.Item(1, 2).Value = 0.09439885811184578
$ws2.Cells.Item(1, 3).Value = 0.9056011418881542
